$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old second header row ("(m3/s)", "(MW)", "(MW)", "(GWh)", "(GWh)", "(GWh)")
# so the three canton rows shift up by one and the sheet becomes a single
# flat table (idx / idx2 / Name / Date Start / Date End / ... / GWh columns).
$ws.Rows.Item(2).Delete()

# Start the header row formatting clean, then rebuild it: two new leading
# index columns plus the renamed / split unit headers.
$ws.Rows.Item(1).ClearFormats()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# The unit-style header cells (m3/s, MW, GWh columns) keep the workbook's
# small-font header look.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Match the recorded selection state after the edit.
$ws.Range("A2:K2").Select()
